$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3, shifting rows 3:25 down to 4:26
$ws.Rows("3:3").Insert()

# Populate the new row 3 with the inserted record's data
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Femacal de La Calera"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44473
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 100112022
$ws.Range("G3").Value = "Arveja Verde"
$ws.Range("H3").Value = "Perfection"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 85
$ws.Range("K3").Value = 35000
$ws.Range("L3").Value = 36000
$ws.Range("M3").Value = 35471
$ws.Range("N3").Value = "$/malla 25 kilos"
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 1419
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"

# Make sure D3 keeps the date style/number format used by the other date cells
$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat
